# "Generate Report for Handback" - refresh the timestamp cells that the
# handback report regenerates on each run.
#
# Sheet "Overview", row 3 (2dc185da-...md):
#   G3 "Latest HO Xliff Generate Date"     2016-08-31 08:55:31 -> 2016-08-31 08:56:26
#
# Sheet "zh-cn", row 3 (2dc185da-...zh-cn.xlf):
#   H3 "Correspond Handoff Datetime"       2016-08-31 08:55:26 -> 2016-08-31 08:56:21
#   K3 "Correspond Handback DateTime"      2016-08-31 08:55:43 -> 2016-08-31 08:56:40
#
# Sheet "de-de", row 3 (2dc185da-...de-de.xlf):
#   H3 "Correspond Handoff Datetime"       2016-08-31 08:55:31 -> 2016-08-31 08:56:47
#   (K3 "Correspond Handback DateTime" is untouched by this run.)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-31 08:56:26"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-31 08:56:21"
$zhcn.Range("K3").Value = "2016-08-31 08:56:40"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-31 08:56:47"
